$wb = $excel.ActiveWorkbook

# --- Overall sheet: A2 number -> text with thousands separator ---
$ws = $wb.Worksheets.Item("Overall")
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "3,981"

# --- County sheet: column B numbers -> text (rows 2-65) ---
$ws = $wb.Worksheets.Item("County")
$ws.Range("B2:B65").NumberFormat = "@"
$ws.Range("B2").Value = "115"
$ws.Range("B3").Value = "4"
$ws.Range("B4").Value = "32"
$ws.Range("B5").Value = "4"
$ws.Range("B6").Value = "100"
$ws.Range("B7").Value = "262"
$ws.Range("B8").Value = "3"
$ws.Range("B9").Value = "32"
$ws.Range("B10").Value = "26"
$ws.Range("B11").Value = "23"
$ws.Range("B12").Value = "75"
$ws.Range("B13").Value = "11"
$ws.Range("B14").Value = "2"
$ws.Range("B15").Value = "2"
$ws.Range("B16").Value = "212"
$ws.Range("B17").Value = "62"
$ws.Range("B18").Value = "17"
$ws.Range("B19").Value = "6"
$ws.Range("B20").Value = "12"
$ws.Range("B21").Value = "1"
$ws.Range("B22").Value = "4"
$ws.Range("B23").Value = "2"
$ws.Range("B24").Value = "7"
$ws.Range("B25").Value = "5"
$ws.Range("B26").Value = "17"
$ws.Range("B27").Value = "20"
$ws.Range("B28").Value = "277"
$ws.Range("B29").Value = "3"
$ws.Range("B30").Value = "45"
$ws.Range("B31").Value = "7"
$ws.Range("B32").Value = "4"
$ws.Range("B33").Value = "46"
$ws.Range("B34").Value = "117"
$ws.Range("B35").Value = "160"
$ws.Range("B36").Value = "8"
$ws.Range("B37").Value = "1"
$ws.Range("B38").Value = "8"
$ws.Range("B39").Value = "63"
$ws.Range("B40").Value = "58"
$ws.Range("B41").Value = "40"
$ws.Range("B42").Value = "513"
$ws.Range("B43").Value = "67"
$ws.Range("B44").Value = "16"
$ws.Range("B45").Value = "40"
$ws.Range("B46").Value = "5"
$ws.Range("B47").Value = "269"
$ws.Range("B48").Value = "38"
$ws.Range("B49").Value = "294"
$ws.Range("B50").Value = "60"
$ws.Range("B51").Value = "247"
$ws.Range("B52").Value = "101"
$ws.Range("B53").Value = "12"
$ws.Range("B54").Value = "8"
$ws.Range("B55").Value = "137"
$ws.Range("B56").Value = "61"
$ws.Range("B57").Value = "44"
$ws.Range("B58").Value = "47"
$ws.Range("B59").Value = "10"
$ws.Range("B60").Value = "6"
$ws.Range("B61").Value = "4"
$ws.Range("B62").Value = "87"
$ws.Range("B63").Value = "3"
$ws.Range("B64").Value = "15"
$ws.Range("B65").Value = "4"

# --- County sheet: row 66 (Glades County) -> formatted text, all zero ---
$ws.Range("B66:F66").NumberFormat = "@"
$ws.Range("B66").Value = "0.00%"
$ws.Range("C66").Value = "$0"
$ws.Range("D66").Value = "0.00%"
$ws.Range("E66").Value = "0.00%"
$ws.Range("F66").Value = "0.00%"

# --- County sheet: new row 67 (Total) ---
$ws.Range("A67:F67").NumberFormat = "@"
$ws.Range("A67").Value = "Total"
$ws.Range("B67").Value = "3,981"
$ws.Range("C67").Value = "$12,417,989,542"
$ws.Range("D67").Value = "7.10%"
$ws.Range("E67").Value = "-16.26%"
$ws.Range("F67").Value = "69.86%"

# --- Congressional District sheet: column B numbers -> text (rows 2-30) ---
$ws = $wb.Worksheets.Item("Congressional District")
$ws.Range("B2:B30").NumberFormat = "@"
$ws.Range("B2").Value = "119"
$ws.Range("B3").Value = "220"
$ws.Range("B4").Value = "75"
$ws.Range("B5").Value = "95"
$ws.Range("B6").Value = "155"
$ws.Range("B7").Value = "266"
$ws.Range("B8").Value = "104"
$ws.Range("B9").Value = "85"
$ws.Range("B10").Value = "177"
$ws.Range("B11").Value = "123"
$ws.Range("B12").Value = "166"
$ws.Range("B13").Value = "257"
$ws.Range("B14").Value = "138"
$ws.Range("B15").Value = "135"
$ws.Range("B16").Value = "139"
$ws.Range("B17").Value = "159"
$ws.Range("B18").Value = "147"
$ws.Range("B19").Value = "71"
$ws.Range("B20").Value = "106"
$ws.Range("B21").Value = "233"
$ws.Range("B22").Value = "110"
$ws.Range("B23").Value = "202"
$ws.Range("B24").Value = "130"
$ws.Range("B25").Value = "164"
$ws.Range("B26").Value = "119"
$ws.Range("B27").Value = "84"
$ws.Range("B28").Value = "145"
$ws.Range("B29").Value = "57"
$ws.Range("B30").Value = "3,981"

# --- Size sheet: column B numbers -> text (rows 2-8) ---
$ws = $wb.Worksheets.Item("Size")
$ws.Range("B2:B8").NumberFormat = "@"
$ws.Range("B2").Value = "1,060"
$ws.Range("B3").Value = "1,242"
$ws.Range("B4").Value = "653"
$ws.Range("B5").Value = "302"
$ws.Range("B6").Value = "567"
$ws.Range("B7").Value = "157"
$ws.Range("B8").Value = "3,981"

# --- Subsector sheet: column B numbers -> text (rows 2-14) ---
$ws = $wb.Worksheets.Item("Subsector")
$ws.Range("B2:B14").NumberFormat = "@"
$ws.Range("B2").Value = "389"
$ws.Range("B3").Value = "591"
$ws.Range("B4").Value = "167"
$ws.Range("B5").Value = "421"
$ws.Range("B6").Value = "24"
$ws.Range("B7").Value = "1,247"
$ws.Range("B8").Value = "50"
$ws.Range("B9").Value = "2"
$ws.Range("B10").Value = "276"
$ws.Range("B11").Value = "95"
$ws.Range("B12").Value = "679"
$ws.Range("B13").Value = "40"
$ws.Range("B14").Value = "3,981"
